$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 73. This shifts the existing data in rows 73:213
# down to rows 74:214 (carrying values AND formatting, exactly matching the
# target diff, since every row below 73 just takes on the prior row's data).
$ws.Rows("73:73").Insert()

# Populate the newly-inserted (blank) row 73 with the new weekly record.
$ws.Cells.Item(73, 1).Value = 5
$ws.Cells.Item(73, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(73, 3).Value = "Maule"
$ws.Cells.Item(73, 4).Value = 44544
$ws.Cells.Item(73, 5).Value = 7
$ws.Cells.Item(73, 6).Value = 100112003
$ws.Cells.Item(73, 7).Value = "Ajo"
$ws.Cells.Item(73, 8).Value = "Chino"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 200
$ws.Cells.Item(73, 11).Value = 20000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 13).Value = 20000
$ws.Cells.Item(73, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(73, 15).Value = "China"
$ws.Cells.Item(73, 16).Value = 2000
$ws.Cells.Item(73, 17).Value = 10
$ws.Cells.Item(73, 18).Value = "Hortaliza"
